$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the Date column (A2:A4) with 11/1/2018, reusing one shared date style
# (numFmtId 14, the built-in short-date format) across all three cells.
$ws.Range("A2").NumberFormat = "mm-dd-yy"
$ws.Range("A2").Value = 43405
$ws.Range("A2").Copy()
$ws.Range("A3:A4").PasteSpecial(-4122)
$ws.Range("A3").Value = 43405
$ws.Range("A4").Value = 43405

# Update the active selection to match the new data entry
$ws.Range("A2:F4").Select()
